# Weekly update: insert 3 new price rows for "Palta" (Hass, Provincia de Quillota)
# at the top of the existing data block (row 736), pushing the older rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows at 736:738 -- existing rows 736..782 shift down to 739..785.
$ws.Rows("736:738").Insert()

# Column map:
# A Mercado ID | B Mercado | C Region | D Fecha | E Codreg | F Tipo
# G Producto ID | H Producto | I Categoria ID | J Categoria | K Variedad
# L Calidad | M Volumen | N Precio minimo | O Precio maximo
# P Precio promedio ponderado | Q Unidad de comercializacion | R Origen
# S Precio $/Kg | T Kg / unidad

$newRows = @(
    @{ Row=736; D=44585; K="Hass"; L="Especial"; M=250; N=3000; O=3000; P=3000; Q="`$/kilo (en caja de 20 kilos)"; R="Provincia de Quillota"; S=3000; T=1 },
    @{ Row=737; D=44585; K="Hass"; L="Primera";  M=300; N=2700; O=2700; P=2700; Q="`$/kilo (en caja de 20 kilos)"; R="Provincia de Quillota"; S=2700; T=1 },
    @{ Row=738; D=44585; K="Hass"; L="Segunda";  M=120; N=2500; O=2500; P=2500; Q="`$/kilo (en caja de 20 kilos)"; R="Provincia de Quillota"; S=2500; T=1 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value2  = 5
    $ws.Cells.Item($row, 2).Value2  = "Macroferia Regional de Talca"
    $ws.Cells.Item($row, 3).Value2  = "Maule"
    $ws.Cells.Item($row, 4).Value2  = $r.D
    $ws.Cells.Item($row, 5).Value2  = 7
    $ws.Cells.Item($row, 6).Value2  = "Fruta"
    $ws.Cells.Item($row, 7).Value2  = 100106
    $ws.Cells.Item($row, 8).Value2  = "Oleaginosos"
    $ws.Cells.Item($row, 9).Value2  = 100106002
    $ws.Cells.Item($row, 10).Value2 = "Palta"
    $ws.Cells.Item($row, 11).Value2 = $r.K
    $ws.Cells.Item($row, 12).Value2 = $r.L
    $ws.Cells.Item($row, 13).Value2 = $r.M
    $ws.Cells.Item($row, 14).Value2 = $r.N
    $ws.Cells.Item($row, 15).Value2 = $r.O
    $ws.Cells.Item($row, 16).Value2 = $r.P
    $ws.Cells.Item($row, 17).Value2 = $r.Q
    $ws.Cells.Item($row, 18).Value2 = $r.R
    $ws.Cells.Item($row, 19).Value2 = $r.S
    $ws.Cells.Item($row, 20).Value2 = $r.T
}
